# Update the "dSF" (column F) values for a handful of rows, per
# "repull data, push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = -1
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -5
$ws.Range("F17").Value = 2
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = -7
